# Performance.xlsx update: add the latest "fedora -> mint" migration
# benchmark rows (92-94) to the "Initial Position Single Thread" sheet,
# mirroring the existing block at rows 88-90, and record a note about the
# OS migration in column P.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Initial Position Single Thread")

$xlPasteFormats = -4122

# ---- carry over formatting from the previous benchmark block ---------
# (rows 88-90 have the same look we want for the new rows 92-94; only
# copy the columns that are actually populated in each source row so we
# don't manufacture empty styled cells where the source has none)
$ws.Range("A88:N88").Copy()
$ws.Range("A92:N92").PasteSpecial($xlPasteFormats)
$ws.Range("C89:N89").Copy()
$ws.Range("C93:N93").PasteSpecial($xlPasteFormats)
$ws.Range("I90:N90").Copy()
$ws.Range("I94:N94").PasteSpecial($xlPasteFormats)

# the "gain" columns (G/M) use a green "good" style when the run improved
# on the previous one and a red "bad" style when it regressed - this
# batch regressed, so pick up the red template already used elsewhere.
$ws.Range("M22").Copy()
$ws.Range("G92").PasteSpecial($xlPasteFormats)
$ws.Range("M22").Copy()
$ws.Range("M92").PasteSpecial($xlPasteFormats)
$ws.Range("M23").Copy()
$ws.Range("G93").PasteSpecial($xlPasteFormats)
$ws.Range("M23").Copy()
$ws.Range("M93").PasteSpecial($xlPasteFormats)

# ---- row 92 (classic depth 4 + bulk depth 4) --------------------------
$ws.Range("A92").Value = 45991
$ws.Range("C92").Value = 4
$ws.Range("D92").Value = 206603
$ws.Range("E92").Value = 207
$ws.Range("F92").Formula = "=D92/E92*1000"
$ws.Range("G92").Formula = "=(E88-E92)/E88"
$ws.Range("H92").Formula = "=(F92-80000000)/80000000"
$ws.Range("I92").Value = 4
$ws.Range("J92").Value = 197281
$ws.Range("K92").Value = 6
$ws.Range("L92").Formula = "=J92/K92*1000"
$ws.Range("M92").Formula = "=(K88-K92)/K88"
$ws.Range("N92").Formula = "=(L92-80000000)/80000000"
$ws.Range("P92").Value = "fedora to debian mint"

# ---- row 93 (classic depth 5 + bulk depth 5) --------------------------
$ws.Range("C93").Value = 5
$ws.Range("D93").Value = 5072212
$ws.Range("E93").Value = 5126
$ws.Range("F93").Formula = "=D93/E93*1000"
$ws.Range("G93").Formula = "=(E89-E93)/E89"
$ws.Range("H93").Formula = "=(F93-80000000)/80000000"
$ws.Range("I93").Value = 5
$ws.Range("J93").Value = 4880523
$ws.Range("K93").Value = 167
$ws.Range("L93").Formula = "=J93/K93*1000"
$ws.Range("M93").Formula = "=(K89-K93)/K89"
$ws.Range("N93").Formula = "=(L93-80000000)/80000000"

# ---- row 94 (bulk depth 6 only) ---------------------------------------
$ws.Range("I94").Value = 6
$ws.Range("J94").Value = 119060324
$ws.Range("K94").Value = 4097
$ws.Range("L94").Formula = "=J94/K94*1000"
$ws.Range("M94").Formula = "=(K90-K94)/K90"
$ws.Range("N94").Formula = "=(L94-80000000)/80000000"

[void]$ws.Range("J94").Select()
